function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6: change the table's style (gallery swap) from "Table_0" to the
#    built-in style {C28D5F3D-3166-44D5-8779-E896A65F9E4C}.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{C28D5F3D-3166-44D5-8779-E896A65F9E4C}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the presentation's (slide-master) theme so its colour scheme
#    matches the standard "Office Theme" palette.
# ---------------------------------------------------------------------------
$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Colors(1).RGB  = (RGB 0x00 0x00 0x00)   # dk1      000000
$cs.Colors(2).RGB  = (RGB 0xFF 0xFF 0xFF)   # lt1      FFFFFF
$cs.Colors(3).RGB  = (RGB 0x44 0x54 0x6A)   # dk2      44546A
$cs.Colors(4).RGB  = (RGB 0xE7 0xE6 0xE6)   # lt2      E7E6E6
$cs.Colors(5).RGB  = (RGB 0x5B 0x9B 0xD5)   # accent1  5B9BD5
$cs.Colors(6).RGB  = (RGB 0xED 0x7D 0x31)   # accent2  ED7D31
$cs.Colors(7).RGB  = (RGB 0xA5 0xA5 0xA5)   # accent3  A5A5A5
$cs.Colors(8).RGB  = (RGB 0xFF 0xC0 0x00)   # accent4  FFC000
$cs.Colors(9).RGB  = (RGB 0x44 0x72 0xC4)   # accent5  4472C4
$cs.Colors(10).RGB = (RGB 0x70 0xAD 0x47)   # accent6  70AD47
$cs.Colors(11).RGB = (RGB 0x05 0x63 0xC1)   # hlink    0563C1
$cs.Colors(12).RGB = (RGB 0x95 0x4F 0x72)   # folHlink 954F72
